# Add Consensus Economics inflation forecasts
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the inflation forecast columns
$ws.Range("L1").Value = "forecast_inf_1step"
$ws.Range("M1").Value = "forecast_inf_2step"
$ws.Range("K1").Value = "forecast_inf_current"

# Fill forecast_date (G) plus GDP (H:J) and new inflation (K:M) forecasts
# for every data row (2-48) with the Consensus Economics November 2020 figures.
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 7).Value = "November 2020"
    $ws.Cells.Item($r, 8).Value = -3.852
    $ws.Cells.Item($r, 9).Value = 3.013
    $ws.Cells.Item($r, 10).Value = 2.334
    $ws.Cells.Item($r, 11).Value = 3.916
    $ws.Cells.Item($r, 12).Value = 3.595
    $ws.Cells.Item($r, 13).Value = 3.576
}

[void]$ws.Range("L10").Select()
